# Updates output/violent-crime-full-year.xlsx with 2023-06-08 data
# (adds one day's worth of incidents to the running 2023 (column J) totals,
# and a couple of corrected 2022 (column I) totals) across the Citywide,
# By Neighborhood, and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 3111
$ws.Range("J3").Value = 3245
$ws.Range("I4").Value = 1760
$ws.Range("J4").Value = 725
$ws.Range("J5").Value = 254
$ws.Range("J6").Value = 3840
$ws.Range("I7").Value = 26206
$ws.Range("J7").Value = 11175

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 49
$ws.Range("J5").Value = 30
$ws.Range("J7").Value = 336
$ws.Range("J8").Value = 714
$ws.Range("J10").Value = 65
$ws.Range("J15").Value = 128
$ws.Range("J19").Value = 354
$ws.Range("J20").Value = 234
$ws.Range("J26").Value = 20
$ws.Range("J27").Value = 64
$ws.Range("J29").Value = 641
$ws.Range("J31").Value = 83
$ws.Range("J33").Value = 476
$ws.Range("J36").Value = 164
$ws.Range("J42").Value = 446
$ws.Range("J43").Value = 104
$ws.Range("J44").Value = 87
$ws.Range("J48").Value = 111
$ws.Range("J53").Value = 109
$ws.Range("J54").Value = 213
$ws.Range("J55").Value = 145
$ws.Range("J57").Value = 50
$ws.Range("J60").Value = 77
$ws.Range("J61").Value = 17
$ws.Range("I63").Value = 215
$ws.Range("J65").Value = 300
$ws.Range("J67").Value = 401
$ws.Range("J73").Value = 99
$ws.Range("J78").Value = 148
$ws.Range("J79").Value = 327
$ws.Range("J85").Value = 512
$ws.Range("J90").Value = 130
$ws.Range("I101").Value = 26206
$ws.Range("J101").Value = 11175

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 126
$ws.Range("J7").Value = 512

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 83
$ws.Range("J6").Value = 131

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 229
$ws.Range("J5").Value = 23
$ws.Range("J6").Value = 210
$ws.Range("J7").Value = 714

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 112
$ws.Range("J3").Value = 100
$ws.Range("J7").Value = 336

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 85
$ws.Range("J3").Value = 166
$ws.Range("J4").Value = 27
$ws.Range("J7").Value = 401

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 21
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 89
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 150
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 476

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 192
$ws.Range("J3").Value = 221
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 641

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 83
$ws.Range("J3").Value = 99
$ws.Range("J7").Value = 354

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 91
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 446

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 95
$ws.Range("J7").Value = 327

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 25
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 46
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J3").Value = 21
$ws.Range("J6").Value = 36

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 25
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 17
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 17
